$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" (summary) sheet: add a new first data row for 2022-Q4,
#        shifting the existing Q3/Q2/Q1 rows down by one (cells are written
#        directly, without Rows.Insert, so no formatting is inherited from
#        neighbouring rows). ---
$summary = $wb.Worksheets.Item(1)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 24
$summary.Range("D2").Value = 4.63

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 11
$summary.Range("D3").Value = 2.53

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 5
$summary.Range("D4").Value = 0.43

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0

# --- 2. Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. before the
#        existing "2022-Q3" tab), and fill it with the quarterly fund-holding data. ---
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Columns B-G hold text (fund codes like "009714" must keep their leading
# zero, and the numeric-looking figures are stored as text in the source
# data), so force a text number format before writing any of their values.
$q4.Range("B1:G25").NumberFormat = "@"

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$rows = @(
    @(0,  "009714", "华安聚优精选混合",               "64.60", "93.63", "3.78", "2.4419", 5),
    @(1,  "012234", "华安聚弘精选混合A",               "19.69", "94.45", "3.52", "0.6931", 6),
    @(2,  "040016", "华安行业轮动混合",                 "6.67", "92.42", "4.16", "0.2775", 4),
    @(3,  "008371", "华安汇智精选混合",                 "6.67", "94.20", "3.55", "0.2368", 5),
    @(4,  "970010", "方正证券金立方一年持有期混合C",    "10.16", "73.92", "2.33", "0.2367", 9),
    @(5,  "040020", "华安升级主题混合A",                 "4.24", "93.79", "4.47", "0.1895", 4),
    @(6,  "014271", "大成北交所两年定开混合A",           "3.24", "68.93", "3.83", "0.1241", 8),
    @(7,  "008290", "华安现代生活混合",                  "2.74", "93.14", "3.73", "0.1022", 3),
    @(8,  "015071", "鑫元专精特新混合A",                 "2.46", "70.69", "3.25", "0.0800", 8),
    @(9,  "012235", "华安聚弘精选混合C",                 "2.07", "94.45", "3.52", "0.0729", 6),
    @(10, "002319", "大成一带一路灵活配置混合",          "1.49", "92.29", "4.81", "0.0717", 10),
    @(11, "014272", "大成北交所两年定开混合C",           "0.77", "68.93", "3.83", "0.0295", 8),
    @(12, "000056", "建信消费升级混合",                  "0.62", "86.11", "2.64", "0.0164", 10),
    @(13, "001226", "中邮稳健添利灵活配置混合",          "0.40", "93.24", "2.72", "0.0109", 9),
    @(14, "011377", "创金合信积极成长股票A",             "0.17", "93.71", "4.83", "0.0082", 10),
    @(15, "015072", "鑫元专精特新混合C",                 "0.25", "70.69", "3.25", "0.0081", 8),
    @(16, "011438", "红塔红土盛昌优选混合A",             "0.13", "92.67", "4.96", "0.0064", 3),
    @(17, "011378", "创金合信积极成长股票C",             "0.12", "93.71", "4.83", "0.0058", 10),
    @(18, "410009", "华富量子生命力混合",                "0.10", "92.12", "4.77", "0.0048", 8),
    @(19, "015633", "中金景气驱动混合A",                 "0.09", "90.15", "3.64", "0.0033", 3),
    @(20, "014976", "华安升级主题混合C",                 "0.05", "93.79", "4.47", "0.0022", 4),
    @(21, "015634", "中金景气驱动混合C",                 "0.06", "90.15", "3.64", "0.0022", 3),
    @(22, "011439", "红塔红土盛昌优选混合C",             "0.04", "92.67", "4.96", "0.0020", 3),
    @(23, "970009", "方正证券金立方一年持有期混合A",     "0.04", "73.92", "2.33", "0.0009", 9)
)

$r = 2
foreach ($row in $rows) {
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("B$r").Value = $row[1]
    $q4.Range("C$r").Value = $row[2]
    $q4.Range("D$r").Value = $row[3]
    $q4.Range("E$r").Value = $row[4]
    $q4.Range("F$r").Value = $row[5]
    $q4.Range("G$r").Value = $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}
